$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 13 de Agosto de 2020 a las 22:42"

# Update country data rows (values refreshed + reordering of Aruba, Monaco, Islas Malvinas
# that moved to new rank positions because their case counts increased)
# Row 4
$ws.Cells.Item(4,2).Value = 5398560
$ws.Cells.Item(4,3).Value = 38258
$ws.Cells.Item(4,4).Value = 2822162
$ws.Cells.Item(4,5).Value = 2406499
$ws.Cells.Item(4,7).Value = 768
$ws.Cells.Item(4,8).Value = 169899

# Row 22
$ws.Cells.Item(22,2).Value = 222269
$ws.Cells.Item(22,3).Value = 1419
$ws.Cells.Item(22,5).Value = 12188

# Row 33
$ws.Cells.Item(33,2).Value = 89822
$ws.Cells.Item(33,3).Value = 1671
$ws.Cells.Item(33,4).Value = 64746
$ws.Cells.Item(33,5).Value = 24425

# Row 60
$ws.Cells.Item(60,2).Value = 37187
$ws.Cells.Item(60,3).Value = 488
$ws.Cells.Item(60,4).Value = 26004
$ws.Cells.Item(60,5).Value = 9842
$ws.Cells.Item(60,7).Value = 8
$ws.Cells.Item(60,8).Value = 1341

# Row 62
$ws.Cells.Item(62,2).Value = 33323
$ws.Cells.Item(62,3).Value = 858
$ws.Cells.Item(62,4).Value = 27213
$ws.Cells.Item(62,5).Value = 5894
$ws.Cells.Item(62,7).Value = 6
$ws.Cells.Item(62,8).Value = 216

# Row 69
$ws.Cells.Item(69,2).Value = 26129
$ws.Cells.Item(69,3).Value = 1072
$ws.Cells.Item(69,4).Value = 8412
$ws.Cells.Item(69,5).Value = 17445
$ws.Cells.Item(69,7).Value = 9
$ws.Cells.Item(69,8).Value = 272

# Row 75
$ws.Cells.Item(75,2).Value = 18308
$ws.Cells.Item(75,3).Value = 45
$ws.Cells.Item(75,5).Value = 2587

# Row 76
$ws.Cells.Item(76,2).Value = 16889
$ws.Cells.Item(76,3).Value = 42
$ws.Cells.Item(76,4).Value = 13522
$ws.Cells.Item(76,5).Value = 3260
$ws.Cells.Item(76,7).Value = 2
$ws.Cells.Item(76,8).Value = 107

# Row 90
$ws.Cells.Item(90,2).Value = 8471
$ws.Cells.Item(90,3).Value = 48
$ws.Cells.Item(90,4).Value = 7784
$ws.Cells.Item(90,5).Value = 636
$ws.Cells.Item(90,7).Value = 1
$ws.Cells.Item(90,8).Value = 51

# Row 91
$ws.Cells.Item(91,2).Value = 8198
$ws.Cells.Item(91,3).Value = 82
$ws.Cells.Item(91,4).Value = 7120
$ws.Cells.Item(91,5).Value = 1028

# Row 121
$ws.Cells.Item(121,2).Value = 3073
$ws.Cells.Item(121,3).Value = 73
$ws.Cells.Item(121,4).Value = 2232
$ws.Cells.Item(121,5).Value = 808

# Row 159 -> Aruba
$ws.Cells.Item(159,1).Value = "Aruba"
$ws.Cells.Item(159,2).Value = 894
$ws.Cells.Item(159,3).Value = 96
$ws.Cells.Item(159,4).Value = 114
$ws.Cells.Item(159,5).Value = 776
$ws.Cells.Item(159,8).Value = 4

# Row 160 -> Lesoto
$ws.Cells.Item(160,1).Value = "Lesoto"
$ws.Cells.Item(160,2).Value = 884
$ws.Cells.Item(160,3).Value = 86
$ws.Cells.Item(160,4).Value = 271
$ws.Cells.Item(160,5).Value = 588
$ws.Cells.Item(160,7).Value = 1
$ws.Cells.Item(160,8).Value = 25

# Row 161 -> Santo Tome y Principe
$ws.Cells.Item(161,1).Value = "Santo Tome y Principe"
$ws.Cells.Item(161,2).Value = 882
$ws.Cells.Item(161,4).Value = 807
$ws.Cells.Item(161,5).Value = 60
$ws.Cells.Item(161,8).Value = 15

# Row 187 -> Monaco
$ws.Cells.Item(187,1).Value = "Monaco"
$ws.Cells.Item(187,3).Value = 3
$ws.Cells.Item(187,4).Value = 114
$ws.Cells.Item(187,5).Value = 26
$ws.Cells.Item(187,8).Value = 4

# Row 188 -> Barbados
$ws.Cells.Item(188,1).Value = "Barbados"
$ws.Cells.Item(188,2).Value = 144
$ws.Cells.Item(188,4).Value = 115
$ws.Cells.Item(188,5).Value = 22
$ws.Cells.Item(188,8).Value = 7

# Row 189 -> Brunei
$ws.Cells.Item(189,1).Value = "Brunei"
$ws.Cells.Item(189,2).Value = 142
$ws.Cells.Item(189,4).Value = 138
$ws.Cells.Item(189,5).Value = 1
$ws.Cells.Item(189,8).Value = 3

# Row 213 -> Islas Malvinas
$ws.Cells.Item(213,1).Value = "Islas Malvinas"
$ws.Cells.Item(213,4).Value = 13
$ws.Cells.Item(213,8).Value = 0

# Row 214 -> Montserrat
$ws.Cells.Item(214,1).Value = "Montserrat"
$ws.Cells.Item(214,4).Value = 12
$ws.Cells.Item(214,8).Value = 1
